# Applies scheduled market-price/profit refresh to the Seraph_Profits workbook.
# For each leve row listed below, updates currentAveragePrice(NQ/HQ) (H/I/J) and
# LevePrice/LeveProfit (K/L/M/N) columns to the latest observed values; cells that
# the refresh leaves blank are cleared instead of written as 0.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 18: Growth Formula Beta
$ws.Range("H18").Value = 937.5
$ws.Range("I18").Value = 937.5
$ws.Range("K18").Value = 937.5
$ws.Range("M18").Value = -653.5

# ALC row 40: Horn Glue
$ws.Range("H40").Value = 1490.4517
$ws.Range("I40").Value = 1492.7142
$ws.Range("K40").Value = 1492.7142
$ws.Range("M40").Value = -1317.7142

# ALC row 42: Hi-Potion of Dexterity
$ws.Range("H42").Value = 300.14285
$ws.Range("I42").Value = 210.6
$ws.Range("J42").Value = 524
$ws.Range("K42").Value = 631.8
$ws.Range("L42").Value = 1572
$ws.Range("M42").Value = -401.8
$ws.Range("N42").Value = -2032

# ALC row 45: Blinding Potion
$ws.Range("H45").Value = 500
$ws.Range("J45").Value = 500
$ws.Range("L45").Value = 1500
$ws.Range("N45").Value = -1884

# ALC row 51: Shark Oil
$ws.Range("H51").Value = 20001
$ws.Range("I51").Value = 20001
$ws.Range("K51").Value = 20001
$ws.Range("M51").Value = -19517

# ALC row 58: Mega-Potion of Vitality
$ws.Range("H58").Value = 5352.5
$ws.Range("I58").Value = 57.5
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 172.5
$ws.Range("L58").Value = 24000
$ws.Range("M58").Value = -22.5
$ws.Range("N58").Value = -24300

# ALC row 99: Commanding Craftsman's Tea
$ws.Range("H99").Value = 492
$ws.Range("I99").Value = 485
$ws.Range("J99").Value = 499
$ws.Range("K99").Value = 1455
$ws.Range("L99").Value = 1497
$ws.Range("M99").Value = 43
$ws.Range("N99").Value = -4493

# ALC row 106: Enchanted Palladium Ink
$ws.Range("H106").Value = 19232.738
$ws.Range("I106").Value = 23354.334
$ws.Range("J106").Value = 4395
$ws.Range("K106").Value = 23354.334
$ws.Range("L106").Value = 4395
$ws.Range("M106").Value = -22723.334
$ws.Range("N106").Value = -5657

# ALC row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2034.1333
$ws.Range("J138").Value = 2089.75
$ws.Range("L138").Value = 6269.25
$ws.Range("N138").Value = -16549.25

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45: Mythril Ingot
$ws.Range("H45").Value = 3111.125
$ws.Range("I45").Value = 1644.5
$ws.Range("K45").Value = 1644.5
$ws.Range("M45").Value = -1267.5

# ARM row 122: High Durium Nugget
$ws.Range("H122").Value = 1273238.4
$ws.Range("I122").Value = 2013997.4
$ws.Range("K122").Value = 6041992.199999999
$ws.Range("M122").Value = -6039542.199999999

# ARM row 126: Bismuth Ingot
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20: Iron Ingot
$ws.Range("H20").Value = 4997
$ws.Range("I20").Value = 4995
$ws.Range("J20").Value = 4999
$ws.Range("K20").Value = 4995
$ws.Range("L20").Value = 4999
$ws.Range("M20").Value = -4748
$ws.Range("N20").Value = -5493

# BSM row 99: Oroshigane Ingot
$ws.Range("H99").Value = 1877.7858
$ws.Range("I99").Value = 2032
$ws.Range("J99").Value = 1492.25
$ws.Range("K99").Value = 2032
$ws.Range("L99").Value = 1492.25
$ws.Range("M99").Value = -534
$ws.Range("N99").Value = -4488.25

# BSM row 105: Molybdenum Ingot
$ws.Range("H105").Value = 2875.2942
$ws.Range("I105").Value = 2499.2856
$ws.Range("K105").Value = 2499.2856
$ws.Range("M105").Value = -752.2856000000002

$ws = $wb.Worksheets.Item("CRP")
# CRP row 3: Maple Pattens
$ws.Range("H3").Value = 20620
$ws.Range("I3").Value = 18992.334
$ws.Range("J3").Value = 21596.6
$ws.Range("K3").Value = 18992.334
$ws.Range("L3").Value = 21596.6
$ws.Range("M3").Value = -18879.334
$ws.Range("N3").Value = -21822.6

$ws = $wb.Worksheets.Item("CUL")
# CUL row 2: Table Salt
$ws.Range("H2").Value = 473.85715
$ws.Range("I2").Value = 47.5
$ws.Range("J2").Value = 644.4
$ws.Range("K2").Value = 285
$ws.Range("L2").Value = 3866.4
$ws.Range("M2").Value = -172
$ws.Range("N2").Value = -4092.4

# CUL row 4: Boiled Egg
$ws.Range("H4").Value = 23300310
$ws.Range("I4").Value = 30152542
$ws.Range("K4").Value = 90457626
$ws.Range("M4").Value = -90457514

# CUL row 9: Jack-o'-lantern
$ws.Range("H9").Value = 236
$ws.Range("I9").Value = 225
$ws.Range("J9").Value = 258
$ws.Range("K9").Value = 675
$ws.Range("L9").Value = 774
$ws.Range("M9").Value = -451
$ws.Range("N9").Value = -1222

# CUL row 12: Kukuru Butter
$ws.Range("H12").Value = 156.6875
$ws.Range("I12").Value = 147.5
$ws.Range("J12").Value = 172
$ws.Range("K12").Value = 442.5
$ws.Range("L12").Value = 516
$ws.Range("M12").Value = -269.5
$ws.Range("N12").Value = -862

# CUL row 36: Crumpet
$ws.Range("H36").Value = 400
$ws.Range("I36").Value = 400
$ws.Range("K36").Value = 1200
$ws.Range("M36").Value = -1031

# CUL row 37: Eel Pie
$ws.Range("H37").Value = 99000
$ws.Range("J37").Value = 99000
$ws.Range("L37").Value = 297000
$ws.Range("N37").Value = -297224

# CUL row 43: Baked Sole
$ws.Range("H43").Value = 2000
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 6000
$ws.Range("N43").Value = -6228

# CUL row 101: Egg Foo Young
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# CUL row 121: Coffee Biscuit
$ws.Range("H121").Value = 558.8
$ws.Range("J121").Value = 1500
$ws.Range("L121").Value = 4500
$ws.Range("N121").Value = -7120

# CUL row 139: Wild Banana Blend
$ws.Range("H139").Value = 1992.7142
$ws.Range("I139").Value = 1992.7142
$ws.Range("K139").Value = 5978.142599999999
$ws.Range("M139").Value = -838.1425999999992

$ws = $wb.Worksheets.Item("GSM")
# GSM row 11: Copper Ring
$ws.Range("H11").Value = 35000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# GSM row 29: Brass Ear Cuffs
$ws.Range("H29").Value = 7250
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 9500
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 9500
$ws.Range("M29").Value = -4710
$ws.Range("N29").Value = -10080

# GSM row 102: Durium Ingot
$ws.Range("H102").Value = 2348.2144
$ws.Range("I102").Value = 1786.875
$ws.Range("J102").Value = 3096.6667
$ws.Range("K102").Value = 1786.875
$ws.Range("L102").Value = 3096.6667
$ws.Range("M102").Value = -164.875
$ws.Range("N102").Value = -6340.6667

# GSM row 122: Ametrine
$ws.Range("H122").Value = 127989.5
$ws.Range("I122").Value = 2783.6
$ws.Range("K122").Value = 8350.799999999999
$ws.Range("M122").Value = -5900.799999999999

$ws = $wb.Worksheets.Item("LTW")
# LTW row 25: Hard Leather Ringbands
$ws.Range("H25").Value = 23335.666
$ws.Range("J25").Value = 23335.666
$ws.Range("L25").Value = 23335.666
$ws.Range("N25").Value = -23795.666

# LTW row 40: Toad Leather
$ws.Range("H40").Value = 1803
$ws.Range("I40").Value = 1803
$ws.Range("K40").Value = 1803
$ws.Range("M40").Value = -1667

# LTW row 46: Boar Leather
$ws.Range("H46").Value = 2950
$ws.Range("I46").Value = 4001
$ws.Range("J46").Value = 2833.2222
$ws.Range("K46").Value = 4001
$ws.Range("L46").Value = 2833.2222
$ws.Range("M46").Value = -3813
$ws.Range("N46").Value = -3209.2222

# LTW row 68: Wyvern Leather
$ws.Range("H68").Value = 4837.5
$ws.Range("I68").Value = 6339.4
$ws.Range("J68").Value = 2334.3333
$ws.Range("K68").Value = 6339.4
$ws.Range("L68").Value = 2334.3333
$ws.Range("M68").Value = -5590.4
$ws.Range("N68").Value = -3832.3333

# LTW row 71: Wyvern Leather
$ws.Range("H71").Value = 4837.5
$ws.Range("I71").Value = 6339.4
$ws.Range("J71").Value = 2334.3333
$ws.Range("K71").Value = 31697
$ws.Range("L71").Value = 11671.6665
$ws.Range("M71").Value = -27953
$ws.Range("N71").Value = -19159.6665

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2: Hempen Underpants
$ws.Range("H2").Value = 846.875

# WVR row 5: Hempen Halfgloves
$ws.Range("H5").Value = 1001500
$ws.Range("I5").Value = 1501000
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 1501000
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = -1500888
$ws.Range("N5").Value = -2724

# WVR row 26: Cotton Dress Shoes
$ws.Range("H26").Value = 1011250
$ws.Range("J26").Value = 1011250
$ws.Range("L26").Value = 1011250
$ws.Range("N26").Value = -1011836

# WVR row 31: Cotton Doublet Vest of Crafting
$ws.Range("H31").Value = 30019
$ws.Range("J31").Value = 30019
$ws.Range("L31").Value = 30019
$ws.Range("N31").Value = -30715

# WVR row 33: Velveteen Wedge Cap of Gathering
$ws.Range("H33").Value = 30624.75
$ws.Range("I33").Value = 37999.5
$ws.Range("J33").Value = 23250
$ws.Range("K33").Value = 37999.5
$ws.Range("L33").Value = 23250
$ws.Range("M33").Value = -37749.5
$ws.Range("N33").Value = -23750

# WVR row 36: Velveteen Wedge Cap of Gathering
$ws.Range("H36").Value = 30624.75
$ws.Range("I36").Value = 37999.5
$ws.Range("J36").Value = 23250
$ws.Range("K36").Value = 37999.5
$ws.Range("L36").Value = 23250
$ws.Range("M36").Value = -37749.5
$ws.Range("N36").Value = -23750

# WVR row 64: Rainbow Ribbon of Healing
$ws.Range("H64").Value = 62179.6
$ws.Range("J64").Value = 62179.6
$ws.Range("L64").Value = 62179.6
$ws.Range("N64").Value = -62675.6

# WVR row 67: Rainbow Ribbon of Healing
$ws.Range("H67").Value = 62179.6
$ws.Range("J67").Value = 62179.6
$ws.Range("L67").Value = 62179.6
$ws.Range("N67").Value = -62675.6
